$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row stays the same (Code, Name, Type, Subtype, IsSubledger, SubledgerType, Active, Description, ParentCode)

# New data rows 2-8. In the target workbook EVERY cell in A1:I8 (including
# the numeric-looking codes in column A/I and the former boolean cells in
# E/G) is stored as literal TEXT, not as a number or boolean. A leading
# apostrophe is the standard Excel way to force a value to be treated as
# text rather than being auto-coerced to a number/boolean, and an
# apostrophe on its own forces an empty TEXT cell (instead of a truly
# blank/number cell) for the columns that have no content.
$data = @(
    @("'1000", "Invalid Type",        "INVALID_TYPE", "Current Asset",     "'false", "'", "'true", "Invalid account type",                "'"),
    @("'2000", "'",                   "LIABILITY",    "Current Liability", "'false", "'", "'true", "Missing name field",                  "'"),
    @("'3000", "First Equity Account","EQUITY",        "Equity",            "'false", "'", "'true", "First equity account with code 3000", "'"),
    @("'3000", "Duplicate Code",      "EQUITY",        "Equity",            "'false", "'", "'true", "Duplicate code 3000",                 "'"),
    @("'4000", "Invalid Parent",      "REVENUE",       "Operating Revenue", "'false", "'", "'true", "Invalid parent code",                 "'9999"),
    @("'5000", "Parent Account",      "EXPENSE",       "Operating Expense", "'false", "'", "'true", "Parent account",                      "'"),
    @("'5100", "Type Mismatch",       "REVENUE",       "Operating Revenue", "'false", "'", "'true", "Type mismatch with parent",           "'5000")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $values[$c]
    }
}
